$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column at N, shifting the old N..P data right to O..Q.
$ws.Columns("N").Insert() | Out-Null

# Restore the stored column width for the new column N as close as the
# COM ColumnWidth (character-unit) API allows.
$ws.Range("N1").EntireColumn.ColumnWidth = 9.85

# Move the selection to S9 and make "Repayment schedule" the active sheet/tab
# (this also clears tabSelected on whichever sheet was previously active).
$ws.Range("S9").Select() | Out-Null
$ws.Activate() | Out-Null
